$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look like plain numbers,
# so Excel keeps them as text (matching the original inlineStr cell type)
# instead of silently converting to a floating point number.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.471.79"
$ws.Range("E2").Value = "  -3.25%  "
$ws.Range("D3").Value = "3.702.33"
$ws.Range("E3").Value = "  -3.51%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "596.64"
$ws.Range("E5").Value = "  -2.05%  "
$ws.Range("D6").Value = "165.68"
$ws.Range("E6").Value = "  -5.15%  "
$ws.Range("D7").Value = "3.700.12"
$ws.Range("E7").Value = "  -3.51%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  -3.20%  "
$ws.Range("D11").Value = "6.21"
$ws.Range("E11").Value = "  -3.82%  "
$ws.Range("E12").Value = "  -3.43%  "
$ws.Range("D13").Value = "37.66"
$ws.Range("E13").Value = "  -5.52%  "
$ws.Range("D14").Value = "0.0000241"
$ws.Range("E14").Value = "  -4.98%  "
$ws.Range("D15").Value = "4.308.13"
$ws.Range("E15").Value = "  -3.56%  "
$ws.Range("D16").Value = "3.692.89"
$ws.Range("E16").Value = "  -3.58%  "
$ws.Range("D17").Value = "67.458.70"
$ws.Range("E17").Value = "  -3.37%  "
$ws.Range("D18").Value = "17.63"
$ws.Range("E18").Value = "  +6.20%  "
$ws.Range("D19").Value = "7.18"
$ws.Range("E19").Value = "  -3.50%  "
$ws.Range("E20").Value = "  -3.27%  "
$ws.Range("D21").Value = "492.01"
$ws.Range("E21").Value = "  -2.50%  "
$ws.Range("E22").Value = "  -3.29%  "
$ws.Range("D23").Value = "0.727"
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("D24").Value = "85.80"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").Value = "  -5.93%  "
$ws.Range("D26").Value = "0.0000139"
$ws.Range("E26").Value = "  -1.77%  "
$ws.Range("D27").Value = "12.21"
$ws.Range("E27").Value = "  -3.13%  "
$ws.Range("D28").Value = "10.11"
$ws.Range("E28").Value = "  -2.73%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "2.94"
$ws.Range("E30").Value = "  -1.27%  "
$ws.Range("E31").Value = "  -6.30%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "7.67"
$ws.Range("E32").Value = "  -3.70%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "31.65"
$ws.Range("E33").Value = "  -1.42%  "
$ws.Range("D34").Value = "3.829.63"
$ws.Range("E34").Value = "  -3.61%  "
$ws.Range("D35").Value = "0.108"
$ws.Range("E35").Value = "  -4.52%  "
$ws.Range("D36").Value = "3.629.38"
$ws.Range("E36").Value = "  -3.67%  "
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").Value = "0.997"
$ws.Range("E38").Value = "  -4.37%  "
$ws.Range("D39").Value = "5.77"
$ws.Range("E39").Value = "  -5.01%  "
$ws.Range("E40").Value = "  -6.72%  "
$ws.Range("D41").Value = "0.323"
$ws.Range("E41").Value = "  -3.63%  "
$ws.Range("D42").Value = "433.72"
$ws.Range("E42").Value = "  -11.27%  "
$ws.Range("D43").Value = "48.67"
$ws.Range("E43").Value = "  -2.08%  "
$ws.Range("E44").Value = "  -5.80%  "
$ws.Range("D45").Value = "2.79"
$ws.Range("E45").Value = "  -6.26%  "
$ws.Range("D46").Value = "8.39"
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("D48").Value = "40.66"
$ws.Range("E48").Value = "  -6.03%  "
$ws.Range("D49").Value = "142.08"
$ws.Range("E49").Value = "  +1.49%  "
$ws.Range("D50").Value = "2.758.16"
$ws.Range("E50").Value = "  -5.77%  "
$ws.Range("D51").Value = "0.0348"
$ws.Range("E51").Value = "  -3.21%  "

# Restore default (Normal) cell style on the cells we reformatted as text,
# so no stray number-format style lingers on them.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
